$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A201").Value = "TAO-USD"
$ws.Range("A202").Value = "IMX-USD"
$ws.Range("A203").Value = "GRT-USD"
